# Commit message: "add big data with max bundle size=2"
# This script rewrites a number of statistic cells on Sheet3, Sheet4, Sheet5
# and Sheet7 of the workbook with new ("big data") values that correspond to
# a run with a bigger dataset / max bundle size of 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet3 : "Statistics" sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

# Row 3 - Number of matched families (columns AI..AN)
$ws3.Range("AI3").Value = 292
$ws3.Range("AJ3").Value = 272
$ws3.Range("AK3").Value = 294
$ws3.Range("AL3").Value = 294
$ws3.Range("AM3").Value = 262
$ws3.Range("AN3").Value = 304

# Row 4 - Number of matched people (columns AI..AN)
$ws3.Range("AI4").Value = 740
$ws3.Range("AJ4").Value = 871
$ws3.Range("AK4").Value = 738
$ws3.Range("AL4").Value = 735
$ws3.Range("AM4").Value = 931
$ws3.Range("AN4").Value = 672

# Row 6 - Size 1 (columns AI..AN)
$ws3.Range("AI6").Value = 19
$ws3.Range("AJ6").Value = 0
$ws3.Range("AK6").Value = 87
$ws3.Range("AL6").Value = 0
$ws3.Range("AM6").Value = 11
$ws3.Range("AN6").Value = 33

# Row 7 - Size 2 (columns AI..AN)
$ws3.Range("AI7").Value = 152
$ws3.Range("AJ7").Value = 33
$ws3.Range("AK7").Value = 103
$ws3.Range("AL7").Value = 165
$ws3.Range("AM7").Value = 0
$ws3.Range("AN7").Value = 214

# Row 8 - Size 3 (columns AI..AN)
$ws3.Range("AI8").Value = 82
$ws3.Range("AJ8").Value = 151
$ws3.Range("AK8").Value = 0
$ws3.Range("AL8").Value = 120
$ws3.Range("AM8").Value = 125
$ws3.Range("AN8").Value = 23

# Row 9 - Size 4 (columns AI..AN)
$ws3.Range("AI9").Value = 24
$ws3.Range("AJ9").Value = 88
$ws3.Range("AK9").Value = 75
$ws3.Range("AL9").Value = 0
$ws3.Range("AM9").Value = 85
$ws3.Range("AN9").Value = 28

# Row 10 - Size 5 (columns AI..AN)
$ws3.Range("AI10").Value = 15
$ws3.Range("AJ10").Value = 0
$ws3.Range("AK10").Value = 29
$ws3.Range("AL10").Value = 9
$ws3.Range("AM10").Value = 41
$ws3.Range("AN10").Value = 6

# "Scarf" column (J) of the "Average" block
$ws3.Range("J14").Value = 286
$ws3.Range("J15").Value = 781
$ws3.Range("J16").Value = 25
$ws3.Range("J17").Value = 111
$ws3.Range("J18").Value = 84
$ws3.Range("J19").Value = 50
$ws3.Range("J20").Value = 17

# "Scarf" column (J) of the "number of families get i games" block
$ws3.Range("J23").Value = 282
$ws3.Range("J24").Value = 718
$ws3.Range("J25").Value = 0

# "Scarf" column (J) of the "number of people get i games" block
$ws3.Range("J32").Value = 513
$ws3.Range("J33").Value = 2087
$ws3.Range("J34").Value = 0

# ---------------------------------------------------------------------------
# Sheet4 : "Best matched" sheet
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")

# "Family" choice counts (column J)
$ws4.Range("J19").Value = 102
$ws4.Range("J20").Value = 286
$ws4.Range("J21").Value = 272
$ws4.Range("J22").Value = 539
$ws4.Range("J23").Value = 506
$ws4.Range("J24").Value = 13

# "People" choice counts (column J)
$ws4.Range("J28").Value = 291
$ws4.Range("J29").Value = 680
$ws4.Range("J30").Value = 871
$ws4.Range("J31").Value = 1527
$ws4.Range("J32").Value = 1253
$ws4.Range("J33").Value = 65

# ---------------------------------------------------------------------------
# Sheet5 : "Num of families get i-th preferred bundle" sheet
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")

$ws5.Range("G2").Value = 161
$ws5.Range("G3").Value = 61
$ws5.Range("G4").Value = 48
$ws5.Range("G5").Value = 156
$ws5.Range("G6").Value = 183
$ws5.Range("G7").Value = 0
$ws5.Range("G8").Value = 23
$ws5.Range("G11").Value = 24
$ws5.Range("G12").Value = 143
$ws5.Range("G13").Value = 0
$ws5.Range("G14").Value = 62
$ws5.Range("G16").Value = 0
$ws5.Range("G17").Value = 115
$ws5.Range("G18").Value = 24
$ws5.Range("G19").Value = 0
$ws5.Range("G20").Value = 0
$ws5.Range("G21").Value = 0
$ws5.Range("G22").Value = 0
$ws5.Range("G23").Value = 0

# ---------------------------------------------------------------------------
# Sheet7 : "Average bundle rank" / "Standard Deviation" sheet
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Sheet7")

$ws7.Range("G2").Value = 6
$ws7.Range("G5").Value = 61.61
